# Apply the updated cryptocurrency price/volume figures (columns D and E)
# for each affected row, as produced by the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.200.77'
$ws.Cells.Item(2, 5).Value = '  +0.47%  '
$ws.Cells.Item(3, 4).Value = '1.856.05'
$ws.Cells.Item(3, 5).Value = '  +0.64%  '
$ws.Cells.Item(4, 4).Value = "'0.9995"
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).Value = "'0.6975"
$ws.Cells.Item(5, 5).Value = '  +0.83%  '
$ws.Cells.Item(6, 4).Value = "'237.12"
$ws.Cells.Item(6, 5).Value = '  -0.10%  '
$ws.Cells.Item(7, 4).Value = "'0.9997"
$ws.Cells.Item(7, 5).Value = '  -0.20%  '
$ws.Cells.Item(8, 4).Value = "'0.07703"
$ws.Cells.Item(8, 5).Value = '  +1.54%  '
$ws.Cells.Item(9, 4).Value = "'0.3044"
$ws.Cells.Item(9, 5).Value = '  +0.35%  '
$ws.Cells.Item(10, 4).Value = "'23.22"
$ws.Cells.Item(10, 5).Value = '  +0.07%  '
$ws.Cells.Item(11, 4).Value = "'0.08189"
$ws.Cells.Item(11, 5).Value = '  +1.13%  '
$ws.Cells.Item(12, 4).Value = '1.845.95'
$ws.Cells.Item(12, 5).Value = '  -0.61%  '
$ws.Cells.Item(13, 4).Value = "'0.7160"
$ws.Cells.Item(13, 5).Value = '  -0.50%  '
$ws.Cells.Item(14, 4).Value = "'5.149"
$ws.Cells.Item(14, 5).Value = '  -0.28%  '
$ws.Cells.Item(15, 4).Value = "'89.14"
$ws.Cells.Item(15, 5).Value = '  +0.30%  '
$ws.Cells.Item(16, 4).Value = '29.200.53'
$ws.Cells.Item(16, 5).Value = '  -0.22%  '
$ws.Cells.Item(17, 4).Value = "'5.746"
$ws.Cells.Item(17, 5).Value = '  -0.19%  '
$ws.Cells.Item(18, 4).Value = "'13.27"
$ws.Cells.Item(18, 5).Value = '  +1.81%  '
$ws.Cells.Item(19, 4).Value = "'0.000007728"
$ws.Cells.Item(19, 5).Value = '  +0.61%  '
$ws.Cells.Item(20, 4).Value = "'237.34"
$ws.Cells.Item(20, 5).Value = '  -1.83%  '
$ws.Cells.Item(21, 4).Value = "'0.9995"
$ws.Cells.Item(21, 5).Value = '  -0.21%  '
$ws.Cells.Item(22, 4).Value = '2.111.27'
$ws.Cells.Item(22, 5).Value = '  -0.50%  '
$ws.Cells.Item(23, 4).Value = "'1.000"
$ws.Cells.Item(23, 5).Value = '  -0.11%  '
$ws.Cells.Item(24, 4).Value = "'7.420"
$ws.Cells.Item(24, 5).Value = '  -2.26%  '
$ws.Cells.Item(25, 4).Value = "'0.1480"
$ws.Cells.Item(25, 5).Value = '  +1.82%  '
$ws.Cells.Item(26, 4).Value = "'162.40"
$ws.Cells.Item(26, 5).Value = '  +0.89%  '
$ws.Cells.Item(27, 4).Value = "'8.992"
$ws.Cells.Item(27, 5).Value = '  +0.30%  '
$ws.Cells.Item(28, 4).Value = "'17.98"
$ws.Cells.Item(28, 5).Value = '  -0.30%  '
$ws.Cells.Item(29, 4).Value = "'2.036"
$ws.Cells.Item(29, 5).Value = '  +5.89%  '
$ws.Cells.Item(30, 4).Value = "'1.417"
$ws.Cells.Item(30, 5).Value = '  +2.70%  '
$ws.Cells.Item(31, 4).Value = "'4.428"
$ws.Cells.Item(31, 5).Value = '  +0.30%  '
$ws.Cells.Item(32, 4).Value = "'1.477"
$ws.Cells.Item(32, 5).Value = '  -0.92%  '
$ws.Cells.Item(33, 4).Value = "'4.012"
$ws.Cells.Item(33, 5).Value = '  -0.70%  '
$ws.Cells.Item(34, 4).Value = "'0.05191"
$ws.Cells.Item(34, 5).Value = '  -0.60%  '
$ws.Cells.Item(35, 4).Value = "'1.162"
$ws.Cells.Item(35, 5).Value = '  -1.62%  '
$ws.Cells.Item(36, 4).Value = "'0.7095"
$ws.Cells.Item(36, 5).Value = '  +0.23%  '
$ws.Cells.Item(37, 4).Value = "'0.9998"
$ws.Cells.Item(37, 5).Value = '  -0.02%  '
$ws.Cells.Item(38, 4).Value = "'2.660"
$ws.Cells.Item(38, 5).Value = '  -0.09%  '
$ws.Cells.Item(39, 4).Value = "'0.01842"
$ws.Cells.Item(39, 5).Value = '  -0.49%  '
$ws.Cells.Item(40, 4).Value = "'2.724"
$ws.Cells.Item(40, 5).Value = '  +1.37%  '
$ws.Cells.Item(41, 4).Value = "'0.9402"
$ws.Cells.Item(41, 5).Value = '  +2.88%  '
$ws.Cells.Item(42, 4).Value = '1.140.44'
$ws.Cells.Item(42, 5).Value = '  +9.09%  '
$ws.Cells.Item(43, 5).Value = '  -0.07%  '
$ws.Cells.Item(44, 4).Value = "'70.87"
$ws.Cells.Item(44, 5).Value = '  +1.85%  '
$ws.Cells.Item(45, 4).Value = "'5.876"
$ws.Cells.Item(45, 5).Value = '  -1.14%  '
$ws.Cells.Item(46, 5).Value = '  -0.06%  '
$ws.Cells.Item(47, 4).Value = "'103.20"
$ws.Cells.Item(48, 4).Value = "'1.792"
$ws.Cells.Item(48, 5).Value = '  +3.07%  '
$ws.Cells.Item(49, 4).Value = '2.008.16'
$ws.Cells.Item(49, 5).Value = '  -0.32%  '
$ws.Cells.Item(50, 4).Value = "'9.146"
$ws.Cells.Item(50, 5).Value = '  -0.76%  '
$ws.Cells.Item(51, 5).Value = '  -3.31%  '
